# Rename the diff-table column headers from the generic "_old"/"_new"
# suffixes to the concrete format-version suffixes "_FV2210"/"_FV2304",
# then (re-)mark the header row + data range as an Excel Table and freeze
# the header row, matching the upstream commit
# "chore: adapt column header formatting to respective input file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into an Excel Table ------------------------
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$lastCol = $used.Column + $used.Columns.Count - 1
$rng = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$lo = $ws.ListObjects.Add(1, $rng, 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
